$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 104: result "Fallo" / profit -1 ---
$ws.Range("G104").Value = "Fallo"
$ws.Range("H104").Value = -1

# --- Append new rows 127 and 128 ---
# Prevent Excel's automatic date inference for the "fecha" text values so
# they are stored as plain text (matching the existing column B cells),
# then reset the style back to the default (no explicit style index).
$dateCells = $ws.Range("B127:B128")
$dateCells.NumberFormat = "@"

$ws.Range("A127").Value = 14552523
$ws.Range("B127").Value = "2025-09-04"
$ws.Range("C127").Value = "Fajing Sun"
$ws.Range("D127").Value = "Omar Jasika"
$ws.Range("E127").Value = "Gana Omar Jasika"
$ws.Range("F127").Value = 2.1

$ws.Range("A128").Value = 14552906
$ws.Range("B128").Value = "2025-09-04"
$ws.Range("C128").Value = "Eliakim Coulibaly"
$ws.Range("D128").Value = "Billy Harris"
$ws.Range("E128").Value = "Gana Eliakim Coulibaly"
$ws.Range("F128").Value = 3.4

$dateCells.Style = "Normal"
